$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$oldGuid = "2b4f149c-ea37-4f1f-ac12-f9cd1fd61d1f"
$newGuid = "8306a4ad-94d5-447d-bdde-ac27af75e6a4"

$oldZhHash = "715ac7f8c2e9232e81c2b3ee19f6ad47a4261152"
$newHash = "8e820b050197e9776b9e6291a7d217cebc6d06b7"

# ---- Overview sheet ----
$wsOverview.Range("A1").EntireColumn.ColumnWidth = 40
$wsOverview.Range("A2").Value2 = "$newGuid.md"
$wsOverview.Range("B2").Value2 = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value2 = "2016-08-27 20:58:31"

# ---- zh-cn sheet ----
$wsZh.Range("A1").EntireColumn.ColumnWidth = 40
$wsZh.Range("A2").Value2 = "$newGuid.md"
$wsZh.Range("G2").Value2 = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value2 = "2016-08-27 20:58:27"
$wsZh.Range("I2").Value2 = ""
$wsZh.Range("J2").Value2 = ""
$wsZh.Range("K2").Value2 = "0001-01-01 00:00:00"

# ---- de-de sheet ----
$wsDe.Range("A1").EntireColumn.ColumnWidth = 40
$wsDe.Range("A2").Value2 = "$newGuid.md"
$wsDe.Range("G2").Value2 = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value2 = "2016-08-27 20:58:31"
$wsDe.Range("I2").Value2 = ""
$wsDe.Range("J2").Value2 = ""
$wsDe.Range("K2").Value2 = "0001-01-01 00:00:00"
